# edit.ps1 - applies the diff described for draft-gandhi-mpls-rfc6374-sr-03.pptx
#
# Shape Left/Top/Width/Height on this COM host are expressed in points and are
# rounded through a single-precision (f32) float before being converted back
# to EMU (truncating, not rounding-to-nearest). A plain "emu/12700.0" value
# can therefore land one EMU below the intended target. EmuToPt() searches a
# small neighbourhood of the naive point value for one that survives the f32
# round-trip and truncation to reproduce the exact target EMU.
function EmuToPt($targetEmu) {
    $base = $targetEmu / 12700.0
    for ($i = -2000; $i -lt 2000; $i++) {
        $candidate = $base + ($i * 0.0000001)
        $f32 = [single]$candidate
        $emu = [int64][math]::Floor([double]$f32 * 12700.0)
        if ($emu -eq $targetEmu) {
            return $candidate
        }
    }
    return $base
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 7: "Return Path TLV" slide - TextBox 6 (the "TLV is mandatory..."
# callout near the Segment List Sub-TLV figure).
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7TextBox6 = $s7.Shapes.Item(6)
$s7tr = $s7TextBox6.TextFrame.TextRange
$s7para1 = $s7tr.Paragraphs(1)
$s7run1 = $s7para1.Characters(1, 83)
$s7run1.Text = "TLV is mandatory when carried in a probe message and if responder does not support, it MUST return Error "
$s7TextBox6.Left = EmuToPt(234416)
$s7TextBox6.Top = EmuToPt(3398621)
$s7TextBox6.Width = EmuToPt(3664394)
$s7TextBox6.Height = EmuToPt(1323439)

# ---------------------------------------------------------------------
# Slide 8: "Block Number TLV for Loss Measurement" slide.
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# Rectangle 2: ASCII-art diagram box - reposition + "Type TBA2" -> "Type = TBA2"
$s8Rect2 = $s8.Shapes.Item(3)
$s8r2tr = $s8Rect2.TextFrame.TextRange
$s8r2para4 = $s8r2tr.Paragraphs(4)
$s8r2para4.Runs(1).Text = "    |  Type = TBA2  |    Length     | Reserved      | Block Number  |"
$s8Rect2.Left = EmuToPt(685800)
$s8Rect2.Top = EmuToPt(958004)
$s8Rect2.Width = EmuToPt(7772400)
$s8Rect2.Height = EmuToPt(1600438)

# Rectangle 6: bullet list - reposition/resize + "TLV is mandatory..." wording
$s8Rect6 = $s8.Shapes.Item(5)
$s8r6tr = $s8Rect6.TextFrame.TextRange
$s8r6para4 = $s8r6tr.Paragraphs(4)
$s8r6para4.Runs(1).Text = "TLV is mandatory when carried in a probe message and if responder does not support, it MUST return "
$s8Rect6.Left = EmuToPt(609600)
$s8Rect6.Top = EmuToPt(2803062)
$s8Rect6.Width = EmuToPt(8229600)
$s8Rect6.Height = EmuToPt(1815882)

# ---------------------------------------------------------------------
# Slide 14: "PM Probes for SR Policy" slide.
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)

# Content Placeholder 2: reposition/resize + shrink body text from 16pt to 14pt
$s14ContentPH = $s14.Shapes.Item(3)
$s14cptr = $s14ContentPH.TextFrame.TextRange
$s14cptr.Font.Size = 14
$s14ContentPH.Left = EmuToPt(838200)
$s14ContentPH.Top = EmuToPt(3961411)
$s14ContentPH.Width = EmuToPt(7620002)
$s14ContentPH.Height = EmuToPt(778651)

# Rectangle 4: ASCII-art diagram box - figure caption wording
$s14Rect4 = $s14.Shapes.Item(4)
$s14r4tr = $s14Rect4.TextFrame.TextRange
$s14r4para18 = $s14r4tr.Paragraphs(18)
$s14r4para18.Runs(1).Text = "   Figure: Example Probe Message Header for an End-to-end SR-MPLS Policy"

# ---------------------------------------------------------------------
# Slide 15: "PM Probes for P2MP SR Policy" slide.
# ---------------------------------------------------------------------
$s15 = $p.Slides.Item(15)

# Rectangle 2: ASCII-art diagram box - figure caption wording
$s15Rect2 = $s15.Shapes.Item(3)
$s15r2tr = $s15Rect2.TextFrame.TextRange
$s15r2para11 = $s15r2tr.Paragraphs(11)
$s15r2para11.Runs(1).Text = "             Figure: Example Probe Query P2MP SR-MPLS Policy"

# Title 1: slide title wording
$s15Title = $s15.Shapes.Item(4)
$s15Title.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "PM Probe Query for P2MP SR Policy"

# Rectangle 4: bullet list - reposition/resize, new first bullet, font size,
# and insertion of a new second bullet re-using the old first-bullet wording.
$s15Rect4 = $s15.Shapes.Item(5)
$s15r4tr = $s15Rect4.TextFrame.TextRange
$s15r4para1 = $s15r4tr.Paragraphs(1)
$null = $s15r4para1.InsertAfter("`rThe Querier root node sends probe query messages using the Replication Segment for the P2MP SR Policy")
$s15r4para1.Runs(1).Text = "Applicable to one-way delay and loss measurement modes for P2MP SR Policy."
$s15r4tr.Font.Size = 14
$s15Rect4.Left = EmuToPt(495300)
$s15Rect4.Top = EmuToPt(3159740)
$s15Rect4.Width = EmuToPt(8229600)
$s15Rect4.Height = EmuToPt(1467005)
